$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# --- Fix the Cases-tab Neo4j query (drop the erroneous/unused `cohort`
#     match + `Cohort` return column that was causing query errors) ---
$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Chihuahua']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value2 = $casesQuery

# --- Row heights settle lower once the extra "Cohort" lines are gone
#     from the query text (rows re-wrap to fewer visual lines) ---
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 216
$ws.Rows.Item(4).RowHeight = 216

# --- Selection moved from B4 back up to B2, scrolled back to the top
#     of the sheet (no more frozen/scrolled topLeftCell) ---
$ws.Range("B2").Select() | Out-Null
